$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Twitter Bootstrap
$ws.Range("A2").Value = "Twitter Bootstrap"
$ws.Range("B2").Value = "3.2.0"
$ws.Range("D2").Value = "http://getbootstrap.com"
$ws.Range("E2").Value = "Copyright 2011-2015 Twitter, Inc. - Licensed under the MIT license"
$ws.Range("F2").Value = "Front end"

# Row 3 - AngularJS
$ws.Range("A3").Value = "AngularJS"
$ws.Range("B3").Value = "1.4.1"
$ws.Range("D3").Value = "https://angularjs.org/"
$ws.Range("E3").Value = "AngularJS v1.4.1 -   (c) 2010-2015 Google, Inc. http://angularjs.org  - License: MIT"
$ws.Range("F3").Value = "Graphs/Charts"

# Row 4 - C3
$ws.Range("A4").Value = "C3"
$ws.Range("B4").Value = "0.4.10"
$ws.Range("D4").Value = "http://c3js.org/"
$ws.Range("F4").Value = "Graphs/Charts"

# Row 5 - NodeJS
$ws.Range("A5").Value = "NodeJS"
$ws.Range("B5").Value = "0.12.4"
$ws.Range("D5").Value = "https://nodejs.org"
$ws.Range("E5").Value = "NPM Install - The Artistic License 2.0 - Copyright (c) 2000-2006, The Perl Foundation."
$ws.Range("F5").Value = "Application Server"

# Row 6 - Express JS
$ws.Range("A6").Value = "Express JS"
$ws.Range("B6").Value = "4.12.4"
$ws.Range("D6").Value = "https://nodejs.org"
$ws.Range("E6").Value = "NPM Install - The Artistic License 2.0 - Copyright (c) 2000-2006, The Perl Foundation."
$ws.Range("F6").Value = "MVC Web Application Framework for NodeJS"

# Row 7 - Docker
$ws.Range("A7").Value = "Docker"
$ws.Range("B7").Value = "1.7.0"
$ws.Range("D7").Value = "https://www.docker.com/"
$ws.Range("F7").Value = "Container"

# Row 8 - Zabbix
$ws.Range("A8").Value = "Zabbix"
$ws.Range("B8").Value = "2.4.5"
$ws.Range("D8").Value = "http://www.zabbix.com/"
$ws.Range("F8").Value = "Continuous Monitoring"

# Row 11 - Git
$ws.Range("A11").Value = "Git"
$ws.Range("B11").Value = "1.9.5"
$ws.Range("C11").Value = "license.terms"
$ws.Range("D11").Value = "https://git-scm.com/"
$ws.Range("E11").Value = "Regents of the University of California, Sun Microsystems, Inc., and other parties.  "
$ws.Range("F11").Value = "Development"

# Row 12 - Jenkins
$ws.Range("A12").Value = "Jenkins"
$ws.Range("B12").Value = "1.617"
$ws.Range("D12").Value = "https://jenkins-ci.org/"
$ws.Range("F12").Value = "Continuous Integrated Builds"

# Selection matches final authored state
$ws.Range("C11").Select()
